# Scheduled-runner style refresh of live market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across all
# eight Leve-profit tables (one per crafting class). Values below were
# computed externally (e.g. by a market-data poller) and are just
# written back into the corresponding cells of each worksheet/table.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 23392666
$ws.Range("J3").Value = 23392666
$ws.Range("L3").Value = 23392666
$ws.Range("N3").Value = -23392894
$ws.Range("H17").Value = 4321.793
$ws.Range("J17").Value = 4321.793
$ws.Range("L17").Value = 12965.379
$ws.Range("N17").Value = -13301.379
$ws.Range("H62").Value = 2483
$ws.Range("I62").Value = 2139.25
$ws.Range("J62").Value = 3399.6667
$ws.Range("K62").Value = 2139.25
$ws.Range("L62").Value = 3399.6667
$ws.Range("M62").Value = -1515.25
$ws.Range("N62").Value = -4647.6667
$ws.Range("H65").Value = 2483
$ws.Range("I65").Value = 2139.25
$ws.Range("J65").Value = 3399.6667
$ws.Range("K65").Value = 10696.25
$ws.Range("L65").Value = 16998.3335
$ws.Range("M65").Value = -7576.25
$ws.Range("N65").Value = -23238.3335
$ws.Range("H76").Value = 3559
$ws.Range("I76").Value = 3486.375
$ws.Range("K76").Value = 3486.375
$ws.Range("M76").Value = -3171.375
$ws.Range("H79").Value = 3559
$ws.Range("I79").Value = 3486.375
$ws.Range("K79").Value = 3486.375
$ws.Range("M79").Value = -2394.375
$ws.Range("H86").Value = 4683224
$ws.Range("I86").Value = 8192016
$ws.Range("J86").Value = 4835
$ws.Range("K86").Value = 8192016
$ws.Range("L86").Value = 4835
$ws.Range("M86").Value = -8190893
$ws.Range("N86").Value = -7081
$ws.Range("H88").Value = 15395768
$ws.Range("I88").Value = 50005748
$ws.Range("K88").Value = 50005748
$ws.Range("M88").Value = -50005342
$ws.Range("H89").Value = 4683224
$ws.Range("I89").Value = 8192016
$ws.Range("J89").Value = 4835
$ws.Range("K89").Value = 40960080
$ws.Range("L89").Value = 24175
$ws.Range("M89").Value = -40954464
$ws.Range("N89").Value = -35407
$ws.Range("H91").Value = 15395768
$ws.Range("I91").Value = 50005748
$ws.Range("K91").Value = 50005748
$ws.Range("M91").Value = -50004344
$ws.Range("H98").Value = 3029.6667
$ws.Range("I98").Value = 3029.6667
$ws.Range("K98").Value = 3029.6667
$ws.Range("M98").Value = -1531.6667
$ws.Range("H102").Value = 23392666
$ws.Range("J102").Value = 23392666
$ws.Range("L102").Value = 23392666
$ws.Range("N102").Value = -23399156
$ws.Range("H122").Value = 3029.6667
$ws.Range("I122").Value = 3029.6667
$ws.Range("K122").Value = 9089.000100000001
$ws.Range("M122").Value = -6639.000100000001
$ws.Range("H129").Value = 1771.1333
$ws.Range("I129").Value = 513
$ws.Range("J129").Value = 2228.6365
$ws.Range("K129").Value = 1539
$ws.Range("L129").Value = 6685.9095
$ws.Range("M129").Value = 3461
$ws.Range("N129").Value = -16685.9095
$ws.Range("H131").Value = 4025.5625
$ws.Range("J131").Value = 13999.5
$ws.Range("L131").Value = 41998.5
$ws.Range("N131").Value = -52078.5
$ws.Range("H136").Value = 113752
$ws.Range("J136").Value = 113752
$ws.Range("L136").Value = 113752
$ws.Range("N136").Value = -123952
$ws.Range("H140").Value = 61093.6
$ws.Range("J140").Value = 58689.75
$ws.Range("L140").Value = 58689.75
$ws.Range("N140").Value = -69049.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 751607.6
$ws.Range("I2").Value = 1093384.6
$ws.Range("J2").Value = 68053.625
$ws.Range("K2").Value = 1093384.6
$ws.Range("L2").Value = 68053.625
$ws.Range("M2").Value = -1093271.6
$ws.Range("N2").Value = -68279.625
$ws.Range("H32").Value = 3655.3403
$ws.Range("I32").Value = 1707.7
$ws.Range("K32").Value = 1707.7
$ws.Range("M32").Value = -1420.7
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H74").Value = 1009.0333
$ws.Range("I74").Value = 391.53845
$ws.Range("K74").Value = 391.53845
$ws.Range("M74").Value = 482.46155
$ws.Range("H77").Value = 1009.0333
$ws.Range("I77").Value = 391.53845
$ws.Range("K77").Value = 1957.69225
$ws.Range("M77").Value = 2410.30775
$ws.Range("H116").Value = 751607.6
$ws.Range("I116").Value = 1093384.6
$ws.Range("J116").Value = 68053.625
$ws.Range("K116").Value = 1093384.6
$ws.Range("L116").Value = 68053.625
$ws.Range("M116").Value = -1091090.6
$ws.Range("N116").Value = -72641.625
$ws.Range("H122").Value = 3575.7693
$ws.Range("I122").Value = 2296.8948
$ws.Range("K122").Value = 6890.6844
$ws.Range("M122").Value = -4440.6844
$ws.Range("H132").Value = 27119.385
$ws.Range("I132").Value = 32805.89
$ws.Range("K132").Value = 98417.67
$ws.Range("M132").Value = -95887.67
$ws.Range("H140").Value = 112874.75
$ws.Range("J140").Value = 112874.75
$ws.Range("L140").Value = 112874.75
$ws.Range("N140").Value = -123234.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 751607.6
$ws.Range("I3").Value = 1093384.6
$ws.Range("J3").Value = 68053.625
$ws.Range("K3").Value = 1093384.6
$ws.Range("L3").Value = 68053.625
$ws.Range("M3").Value = -1093270.6
$ws.Range("N3").Value = -68281.625
$ws.Range("H20").Value = 4291.5
$ws.Range("I20").Value = 3720.5
$ws.Range("K20").Value = 3720.5
$ws.Range("M20").Value = -3473.5
$ws.Range("H86").Value = 4416.5
$ws.Range("I86").Value = 5500
$ws.Range("J86").Value = 3333
$ws.Range("K86").Value = 5500
$ws.Range("L86").Value = 3333
$ws.Range("M86").Value = -4377
$ws.Range("N86").Value = -5579
$ws.Range("H89").Value = 4416.5
$ws.Range("I89").Value = 5500
$ws.Range("J89").Value = 3333
$ws.Range("K89").Value = 27500
$ws.Range("L89").Value = 16665
$ws.Range("M89").Value = -21884
$ws.Range("N89").Value = -27897
$ws.Range("H105").Value = 2302.889
$ws.Range("I105").Value = 1893.5
$ws.Range("K105").Value = 1893.5
$ws.Range("M105").Value = -146.5
$ws.Range("H134").Value = 2717.973
$ws.Range("I134").Value = 2222.3333
$ws.Range("K134").Value = 6666.999899999999
$ws.Range("M134").Value = -4131.999899999999
$ws.Range("H140").Value = 97978
$ws.Range("J140").Value = 97978
$ws.Range("L140").Value = 97978
$ws.Range("N140").Value = -108338
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1908.7059
$ws.Range("I58").Value = 1986.9
$ws.Range("J58").Value = 1797
$ws.Range("K58").Value = 1986.9
$ws.Range("L58").Value = 1797
$ws.Range("M58").Value = -1783.9
$ws.Range("N58").Value = -2203
$ws.Range("H132").Value = 10107820
$ws.Range("I132").Value = 11497812
$ws.Range("K132").Value = 34493436
$ws.Range("M132").Value = -34490906
$ws.Range("H134").Value = 2114.568
$ws.Range("I134").Value = 2067.6428
$ws.Range("K134").Value = 6202.928400000001
$ws.Range("M134").Value = -3667.928400000001
$ws.Range("H136").Value = 1908.7059
$ws.Range("I136").Value = 1986.9
$ws.Range("J136").Value = 1797
$ws.Range("K136").Value = 5960.700000000001
$ws.Range("L136").Value = 5391
$ws.Range("M136").Value = -3410.700000000001
$ws.Range("N136").Value = -10491
$ws.Range("H141").Value = 91707.91
$ws.Range("J141").Value = 94318.39999999999
$ws.Range("L141").Value = 94318.39999999999
$ws.Range("N141").Value = -104678.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 682.5
$ws.Range("I12").Value = 306.83334
$ws.Range("J12").Value = 843.5
$ws.Range("K12").Value = 920.5000200000001
$ws.Range("L12").Value = 2530.5
$ws.Range("M12").Value = -747.5000200000001
$ws.Range("N12").Value = -2876.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 1101988.1
$ws.Range("I122").Value = 2201862.5
$ws.Range("K122").Value = 6605587.5
$ws.Range("M122").Value = -6603137.5
$ws.Range("H132").Value = 4346.067
$ws.Range("I132").Value = 3388.818
$ws.Range("K132").Value = 10166.454
$ws.Range("M132").Value = -7636.454000000002
$ws.Range("H140").Value = 68439
$ws.Range("J140").Value = 68439
$ws.Range("L140").Value = 68439
$ws.Range("N140").Value = -78799
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1085076.6
$ws.Range("I68").Value = 1750158.2
$ws.Range("K68").Value = 1750158.2
$ws.Range("M68").Value = -1749409.2
$ws.Range("H71").Value = 1085076.6
$ws.Range("I71").Value = 1750158.2
$ws.Range("K71").Value = 8750791
$ws.Range("M71").Value = -8747047
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 199994
$ws.Range("J92").Value = 199994
$ws.Range("L92").Value = 199994
$ws.Range("N92").Value = -204986
$ws.Range("H130").Value = 84499.5
$ws.Range("J130").Value = 84499.5
$ws.Range("L130").Value = 84499.5
$ws.Range("N130").Value = -94539.5
$ws.Range("H132").Value = 19235718
$ws.Range("I132").Value = 2032.8096
$ws.Range("J132").Value = 100017200
$ws.Range("K132").Value = 6098.4288
$ws.Range("L132").Value = 300051600
$ws.Range("M132").Value = -3568.4288
$ws.Range("N132").Value = -300056660
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280
